# Polish Day 1 presentations
# Applies three small text/formatting tweaks:
#   1. Slide 2 - bold the word "implemented" in the cgroups bullet.
#   2. Slide 4 - split "Customer Data Center" label into "Client" / " Data Center".
#   3. Slide 8 - bold a few key words/phrases in the bullet list.

$p = $ppt.ActivePresentation

function Bold-Phrase {
    param($TextRange, $ParagraphIndex, $Phrase)
    $para = $TextRange.Paragraphs($ParagraphIndex)
    $localIdx = $para.Text.IndexOf($Phrase)
    if ($localIdx -lt 0) { return }
    $sub = $TextRange.Characters($para.Start + $localIdx, $Phrase.Length)
    $sub.Font.Bold = $true
}

# ---------------------------------------------------------------------------
# 1. Slide 2, "Linux Containers (LXC) details" content placeholder (shape id 12)
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shp2 = $null
foreach ($sh in $s2.Shapes) {
    if ($sh.Id -eq 12) { $shp2 = $sh }
}
$tr2 = $shp2.TextFrame.TextRange
for ($i = 1; $i -le $tr2.Paragraphs().Count; $i++) {
    if ($tr2.Paragraphs($i).Text -like "Resource management implemented*") {
        Bold-Phrase $tr2 $i "implemented"
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Slide 4, "Rectangle 19" label (shape id 173) - "Customer Data Center"
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $null
foreach ($sh in $s4.Shapes) {
    if ($sh.Id -eq 173) { $shp4 = $sh }
}
$tr4 = $shp4.TextFrame.TextRange
$tr4.Text = "Client" + [char]13 + " Data Center"

# ---------------------------------------------------------------------------
# 3. Slide 8, "Text Placeholder 1" bullet list (shape id 2)
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$shp8 = $null
foreach ($sh in $s8.Shapes) {
    if ($sh.Id -eq 2) { $shp8 = $sh }
}
$tr8 = $shp8.TextFrame.TextRange
for ($i = 1; $i -le $tr8.Paragraphs().Count; $i++) {
    $ptext = $tr8.Paragraphs($i).Text
    if ($ptext -like "Can run on many different platforms*") {
        Bold-Phrase $tr8 $i "different platforms"
    }
    elseif ($ptext -like "Processes share OS resources*") {
        Bold-Phrase $tr8 $i "share"
    }
    elseif ($ptext -like "Isolate the different requirements*") {
        Bold-Phrase $tr8 $i "Isolate"
    }
}
